{"js": "// Fixed #295 Add the version of M2Doc in the template custom properties.\n// The stack-trace text embedded in the document references M2Doc source\n// line numbers that shifted because of this change. Update each\n// \"at ...(File.java:LINE)\" occurrence to its new line number.\n\nconst replacements = [\n  [\"M2DocEvaluator.caseConditional(M2DocEvaluator.java:1267)\",\n   \"M2DocEvaluator.caseConditional(M2DocEvaluator.java:1313)\"],\n  [\"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n   \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1084)\"],\n  [\"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\",\n   \"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1300)\"],\n  [\"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\",\n   \"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:278)\"],\n  [\"M2DocEvaluator.generate(M2DocEvaluator.java:264)\",\n   \"M2DocEvaluator.generate(M2DocEvaluator.java:267)\"],\n  [\"M2DocUtils.generate(M2DocUtils.java:712)\",\n   \"M2DocUtils.generate(M2DocUtils.java:694)\"],\n  [\"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\",\n   \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)\"],\n  [\"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\",\n   \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:384)\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fixed #295 Add the version of M2Doc in the template custom properties.\n# The stack-trace text embedded in the document references M2Doc source\n# line numbers that shifted because of this change. Update each\n# \"at ...(File.java:LINE)\" occurrence to its new line number.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"M2DocEvaluator.caseConditional(M2DocEvaluator.java:1267)\", \"M2DocEvaluator.caseConditional(M2DocEvaluator.java:1313)\"),\n    @(\"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\", \"M2DocEvaluator.doSwitch(M2DocEvaluator.java:1084)\"),\n    @(\"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\", \"M2DocEvaluator.caseBlock(M2DocEvaluator.java:1300)\"),\n    @(\"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\", \"M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:278)\"),\n    @(\"M2DocEvaluator.generate(M2DocEvaluator.java:264)\", \"M2DocEvaluator.generate(M2DocEvaluator.java:267)\"),\n    @(\"M2DocUtils.generate(M2DocUtils.java:712)\", \"M2DocUtils.generate(M2DocUtils.java:694)\"),\n    @(\"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\", \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)\"),\n    @(\"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\", \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:384)\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
